# Reorders data rows 2-13 on the active worksheet (Fruta / Granada sheet).
# Row contents for columns A..T are kept intact; only the row positions that
# each full record occupies are permuted (a reshuffle of the same weekly
# price observations), per mapping: new row r <- old row map[r]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (1-based worksheet rows)
$map = @{
    2  = 3
    3  = 4
    4  = 10
    5  = 8
    6  = 2
    7  = 7
    8  = 5
    9  = 11
    10 = 13
    11 = 12
    12 = 6
    13 = 9
}

$firstCol = 1   # A
$lastCol  = 20  # T

# Snapshot every source row's values (A:T) before any writes, so that
# overwritten cells do not corrupt rows still to be read.
# Value2() is used (as a method call) so date-formatted cells come back as
# raw serial numbers instead of formatted date strings.
$snapshot = @{}
foreach ($srcRow in $map.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $rowVals[$c] = $ws.Cells.Item($srcRow, $c).Value2()
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# Write the permuted rows back using the captured snapshot.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $rowVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $rowVals[$c]
    }
}
